# "Generate Report for Handback"
#
# For the two source files (0986ad25...md and ffffdaf03d5b...md) that were
# previously only "Ready for handoff", the localization round-trip has now
# completed: they've been handed back and are in sync with en-US. This
# script records that on the Overview sheet and on each per-locale sheet
# (zh-cn / de-de), and fills in the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns for the two locale sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdName = "0986ad25-7ee2-4c21-9928-79827311773c.md"
$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/353f1857a989a7f57aa3d6a899b05dc3a1601ed7/e2e/$mdName"

# ---- Overview sheet: update the Status-like columns for both files ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- Per-locale sheets ----
# locale => sheet name; xlf file name; URL of the handed-back xlf; handback datetime
$locales = @(
    @{
        Sheet    = "zh-cn"
        XlfName  = "0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.zh-cn.xlf"
        XlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19c55e992ad8c9e981ed48a84b0e04fada166307/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.zh-cn.xlf"
        Handback = "2016-03-08 12:48:12"
    },
    @{
        Sheet    = "de-de"
        XlfName  = "0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.de-de.xlf"
        XlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d6475ef8e81c04b0b479da5ed20ab4a4341cb59/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.de-de.xlf"
        Handback = "2016-03-08 12:48:23"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Row 2 and row 3 both correspond to the same source/target file pair
    # (0986ad25...md) for this sample data set, and both just got handed
    # back from the same xlf.
    foreach ($row in 2, 3) {
        # Status
        $ws.Cells.Item($row, 2).Value = $newStatus

        # Latest Target File (column E) - the file that was localized
        $ws.Hyperlinks.Add(
            $ws.Cells.Item($row, 5),
            $mdUrl,
            "",
            "",
            $mdName
        ) | Out-Null

        # Latest Handback File (column F) - the xlf handed back from translation
        $ws.Hyperlinks.Add(
            $ws.Cells.Item($row, 6),
            $loc.XlfUrl,
            "",
            "",
            $loc.XlfName
        ) | Out-Null

        # Latest Handback DateTime (column G)
        $ws.Cells.Item($row, 7).Value = $loc.Handback
    }
}
